# Apply weekly fruit/vegetable price update: reorder date/volume/price
# values across rows 2-14 as per the new data snapshot.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").Value = 44215

# Row 3
$ws.Range("D3").Value = 44231
$ws.Range("J3").Value = 250

# Row 4
$ws.Range("D4").Value = 44187
$ws.Range("J4").Value = 160

# Row 5
$ws.Range("D5").Value = 44186
$ws.Range("J5").Value = 160

# Row 6
$ws.Range("D6").Value = 44208
$ws.Range("J6").Value = 160

# Row 7
$ws.Range("D7").Value = 44230
$ws.Range("J7").Value = 250

# Row 8
$ws.Range("D8").Value = 44232
$ws.Range("J8").Value = 250

# Row 9
$ws.Range("D9").Value = 44204
$ws.Range("J9").Value = 430

# Row 11
$ws.Range("D11").Value = 44251
$ws.Range("J11").Value = 120
$ws.Range("L11").Value = 5000
$ws.Range("M11").Value = 5000
$ws.Range("O11").Value = "Región Metropolitana"
$ws.Range("P11").Value = 312

# Row 12
$ws.Range("D12").Value = 44210
$ws.Range("J12").Value = 340

# Row 13
$ws.Range("D13").Value = 44188
$ws.Range("J13").Value = 210
$ws.Range("L13").Value = 6000
$ws.Range("M13").Value = 5500
$ws.Range("O13").Value = "Provincia de Quillota"
$ws.Range("P13").Value = 344

# Row 14
$ws.Range("D14").Value = 44189
